$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 6

$ws.Range("F4").Select()
